$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("C2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.230.92"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.20"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6997"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.75"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08100"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +3.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3018"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.44"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08176"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.859.66"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7061"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.70"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.264.49"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.821"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007905"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.116.36"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.431"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.83"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.879"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1421"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.916"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.359"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.030"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05176"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.160"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7203"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9982"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.688"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01850"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.717"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9356"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.148.61"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +3.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.986"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4251"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.22"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.84"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5292"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.751"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.011.82"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.149"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -0.03%  "
